$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------- Row 4 ----------
$ws.Range("A4").Value = 82950
$ws.Range("P4").Value = "Bränsle, SO om (4), Öl"
$ws.Range("S4").Value = 10
$ws.Range("X4").Value = "Hö-Bor-1968"
$ws.Range("AC4").Value = "CA 08319"
$ws.Range("AR4").Value = ""
$ws.Range("AW4").Value = "Öland- Floraväktarna"
$ws.Range("AX4").Value = "Crister Albinsson"
$ws.Range("AY4").Value = "Floraväkteri Sverige"

# ---------- Row 7 ----------
$ws.Range("A7").Value = 218918
$ws.Range("K7").Value = ""
$ws.Range("P7").Value = "Bränsle, SO om (2), Öl"
$ws.Range("S7").Value = 25
$ws.Range("X7").Value = "Hö-Bor-3930"
$ws.Range("AI7").Value = "vägkant"
$ws.Range("AW7").Value = "Öland- Floraväktarna"
$ws.Range("AX7").Value = "Crister Albinsson"
$ws.Range("AY7").Value = "Floraväkteri Sverige"

# ---------- Row 15 ----------
$ws.Range("A15").Value = 67088883
$ws.Range("I15").Value = ""
$ws.Range("J15").Value = ""
$ws.Range("P15").Value = "Bränsle 500 m SO, Öl"
$ws.Range("S15").Value = 25
$ws.Range("Y15").Value = "'2017-08-09"
$ws.Range("AA15").Value = "'2017-08-09"
$ws.Range("AI15").Value = ""
$ws.Range("AW15").Value = "Stefan Kasselstrand"
$ws.Range("AX15").Value = "Stefan Kasselstrand"

# ---------- Row 20 ----------
$ws.Range("A20").Value = 103312351
$ws.Range("C20").Value = "Godkänd. Foto (eller ljud) granskat av validerare"
$ws.Range("I20").Value = "'1"
$ws.Range("J20").Value = "m²"
$ws.Range("K20").Value = "i frukt"
$ws.Range("P20").Value = "Bränsle 500 m SO, Öl"
$ws.Range("S20").Value = 10
$ws.Range("Y20").Value = "'2022-09-01"
$ws.Range("AA20").Value = "'2022-09-01"
$ws.Range("AC20").Value = ""
$ws.Range("AI20").Value = ""
$ws.Range("AW20").Value = "Stefan Kasselstrand"
$ws.Range("AX20").Value = "Stefan Kasselstrand, Ingvor Kasselstrand, Magnus Kasselstrand"

# ---------- Row 23 ----------
$ws.Range("A23").Value = 95764169
$ws.Range("C23").Value = "Ovaliderad"
$ws.Range("K23").Value = ""
$ws.Range("P23").Value = "Bränsle SO, Öl"
$ws.Range("S23").Value = 25
$ws.Range("X23").Value = "Hö-Bor-8992"
$ws.Range("AI23").Value = ""
$ws.Range("AW23").Value = "Thomas Gunnarsson"
$ws.Range("AX23").Value = "Vera Wendt"
$ws.Range("AY23").Value = "Floraväkteri Sverige"

# ---------- Row 25 ----------
$ws.Range("B25").Value = 103716
